$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Yaya Hotel & Apartments block (cell B11): update Trip Coins amount and "last booked" hours ---
$b11 = $ws.Range("B11").Value()
$b11 = $b11.Replace("Earn ₹ 961.74 in Trip Coins", "Earn ₹ 962.30 in Trip Coins")
$b11 = $b11.Replace("Last booked 5 hrs ago", "Last booked 6 hrs ago")
$ws.Range("B11").Value = $b11

# --- Lux Suites Riara One Residency Angama block (cell B8): update Trip Coins amount ---
$b8 = $ws.Range("B8").Value()
$b8 = $b8.Replace("Earn ₹ 1,322.62 in Trip Coins", "Earn ₹ 1,323.05 in Trip Coins")
$ws.Range("B8").Value = $b8

# --- Kester International Apartment Hotel block (cell A9): update Trip Coins amount ---
$a9 = $ws.Range("A9").Value()
$a9 = $a9.Replace("Earn ₹ 546.73 in Trip Coins", "Earn ₹ 546.53 in Trip Coins")
$ws.Range("A9").Value = $a9

# --- Swap A4 and A5 entries (Woodmere Serviced Apartment <-> Eldon Apartments & Suites) ---
$a4 = $ws.Range("A4").Value()
$a5 = $ws.Range("A5").Value()
$ws.Range("A4").Value = $a5
$ws.Range("A5").Value = $a4
